# Fix bug in RNN_test in 2b
# Updates the "Part 2b" results table (rows 13-14) on Sheet1: the
# SequenceLengthOfTest=100 entry is corrected to SequenceLengthOfTest=400,
# its test-accuracy result is corrected, and the "Custom 1" row is filled
# in with a new experiment (SequenceLengthOfTrain=400) plus commentary.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 13: correct the description + test accuracy -----------------
$ws.Range("C13").Value = "ADAM optimizer with LR=0.001, BatchSize=200, VocabularySize=8000, HiddenUnits=500, SequenceLengthOfTrain=100, SequenceLengthOfTest=400"
$ws.Range("G13").Value = 0.8769

# --- Row 14 ("Custom 1"): new experiment results ----------------------
$ws.Range("C14").Value = "ADAM optimizer with LR=0.001, BatchSize=200, VocabularySize=8000, HiddenUnits=500, SequenceLengthOfTrain=400, SequenceLengthOfTest=400"
$ws.Range("D14").Value = 20
$ws.Range("E14").Value = 0.0175
$ws.Range("F14").Value = 0.9943
$ws.Range("G14").Value = 0.7907
$ws.Range("H14").Value = "The sequence length for training is increased from 100 to 400. The model is now overfitting. Training accuracy is approaching 100% while test accuracy is only 79.07%, which is worse than the given model."
$ws.Range("F14:G14").NumberFormat = "0.00%"
$ws.Rows.Item(14).RowHeight = 65

# --- View state: scroll/selection as left by the author ---------------
$ws.Range("H14").Select()
$excel.ActiveWindow.ScrollRow = 5
